$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (name unchanged)
$ws.Range("B2").Value = 0.8809126784011837
$ws.Range("C2").Value = 0.8809126784011837
$ws.Range("D2").Value = 0.8809126784011837

# Row 3 - RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 0.9953458024937243
$ws.Range("C3").Value = 0.9953147547573437
$ws.Range("D3").Value = 0.9894989695990718

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9892667809184853
$ws.Range("C4").Value = 0.9897283990674534
$ws.Range("D4").Value = 0.9661132147163253

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9982001395604766
$ws.Range("C5").Value = 0.9981167489816167
$ws.Range("D5").Value = 0.9975461422772799
